$wb = $excel.ActiveWorkbook

# The last sheet (Slovakia) is the template for the new "Italy" sheet.
$template = $wb.Worksheets.Item($wb.Worksheets.Count)

# Select the whole sheet before duplicating, mirroring the "select-all"
# artifact left behind on the template sheet once it is no longer active.
$template.Activate() | Out-Null
$template.Cells.Select() | Out-Null

# Duplicate the Slovakia sheet (keeps formatting, merges, styles, etc.)
# placing the copy right after it, then rename it to "Italy".
$template.Copy($null, $template) | Out-Null
$newSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet.Name = "Italy"

# Fill in the Italy-specific market name and Jira/test-case reference.
$newSheet.Range("B2").Value = "Italy Market"
$newSheet.Range("B4").Value = "NGC-3145/T2219/T2221/T2223/T2224 "

# Make the new Italy sheet the active tab.
$newSheet.Activate() | Out-Null
$newSheet.Range("B7").Select() | Out-Null
